# Update leveling-profit sheets with refreshed Universalis price snapshots (scheduled runner sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 146.66667
$ws.Range("I52").Value = 146.66667
$ws.Range("K52").Value = 440.00001
$ws.Range("M52").Value = -280.00001

$ws.Range("H113").Value = 126785.625
$ws.Range("I113").Value = 201657
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 201657
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -198403
$ws.Range("N113").Value = -8508

$ws.Range("H125").Value = 2604.1365
$ws.Range("I125").Value = 2223.9443
$ws.Range("J125").Value = 4315
$ws.Range("K125").Value = 20015.4987
$ws.Range("L125").Value = 38835
$ws.Range("M125").Value = -17555.4987
$ws.Range("N125").Value = -43755

$ws.Range("H132").Value = 6950630
$ws.Range("I132").Value = 7582117
$ws.Range("J132").Value = 4270.3335
$ws.Range("K132").Value = 22746351
$ws.Range("L132").Value = 12811.0005
$ws.Range("M132").Value = -22743821
$ws.Range("N132").Value = -17871.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24714.24
$ws.Range("I32").Value = 4500.4155
$ws.Range("K32").Value = 4500.4155
$ws.Range("M32").Value = -4213.4155

$ws.Range("H74").Value = 1880.6
$ws.Range("I74").Value = 893.2308
$ws.Range("K74").Value = 893.2308
$ws.Range("M74").Value = -19.23080000000004

$ws.Range("H77").Value = 1880.6
$ws.Range("I77").Value = 893.2308
$ws.Range("K77").Value = 4466.154
$ws.Range("M77").Value = -98.15400000000045

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H132").Value = 1529.5306
$ws.Range("I132").Value = 930.7027
$ws.Range("J132").Value = 3375.9167
$ws.Range("K132").Value = 2792.1081
$ws.Range("L132").Value = 10127.7501
$ws.Range("M132").Value = -262.1081000000004
$ws.Range("N132").Value = -15187.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1938.3889
$ws.Range("I134").Value = 1867.9778
$ws.Range("J134").Value = 2290.4443
$ws.Range("K134").Value = 5603.9334
$ws.Range("L134").Value = 6871.3329
$ws.Range("M134").Value = -3068.9334
$ws.Range("N134").Value = -11941.3329

$ws.Range("H135").Value = 36250
$ws.Range("J135").Value = 36250
$ws.Range("L135").Value = 36250
$ws.Range("N135").Value = -46390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27977.334
$ws.Range("I31").Value = 1090.3103
$ws.Range("J31").Value = 53129.71
$ws.Range("K31").Value = 1090.3103
$ws.Range("L31").Value = 53129.71
$ws.Range("M31").Value = -795.3103000000001
$ws.Range("N31").Value = -53719.71

$ws.Range("H34").Value = 27977.334
$ws.Range("I34").Value = 1090.3103
$ws.Range("J34").Value = 53129.71
$ws.Range("K34").Value = 1090.3103
$ws.Range("L34").Value = 53129.71
$ws.Range("M34").Value = -888.3103000000001
$ws.Range("N34").Value = -53533.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 551.5
$ws.Range("I44").Value = 402
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 1206
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -808
$ws.Range("N44").Value = -3796

$ws.Range("H47").Value = 150.4
$ws.Range("I47").Value = 117.166664
$ws.Range("K47").Value = 351.499992
$ws.Range("M47").Value = 79.50000799999998

$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868

$ws.Range("H123").Value = 3986.818
$ws.Range("I123").Value = 2452.5
$ws.Range("J123").Value = 4863.5713
$ws.Range("K123").Value = 7357.5
$ws.Range("L123").Value = 14590.7139
$ws.Range("M123").Value = -4907.5
$ws.Range("N123").Value = -19490.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10979.556
$ws.Range("J52").Value = 10979.556
$ws.Range("L52").Value = 10979.556
$ws.Range("N52").Value = -11497.556

$ws.Range("H70").Value = 100936.29
$ws.Range("I70").Value = 226504.67
$ws.Range("J70").Value = 6760
$ws.Range("K70").Value = 226504.67
$ws.Range("L70").Value = 6760
$ws.Range("M70").Value = -226234.67
$ws.Range("N70").Value = -7300

$ws.Range("H73").Value = 100936.29
$ws.Range("I73").Value = 226504.67
$ws.Range("J73").Value = 6760
$ws.Range("K73").Value = 226504.67
$ws.Range("L73").Value = 6760
$ws.Range("M73").Value = -225568.67
$ws.Range("N73").Value = -8632

$ws.Range("H102").Value = 601865.8
$ws.Range("I102").Value = 2178
$ws.Range("J102").Value = 3000617
$ws.Range("K102").Value = 2178
$ws.Range("L102").Value = 3000617
$ws.Range("M102").Value = -556
$ws.Range("N102").Value = -3003861

$ws.Range("H134").Value = 26806.5
$ws.Range("J134").Value = 26806.5
$ws.Range("L134").Value = 80419.5
$ws.Range("N134").Value = -85489.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7401
$ws.Range("I68").Value = 6001
$ws.Range("J68").Value = 8334.333000000001
$ws.Range("K68").Value = 6001
$ws.Range("L68").Value = 8334.333000000001
$ws.Range("M68").Value = -5252
$ws.Range("N68").Value = -9832.333000000001

$ws.Range("H71").Value = 7401
$ws.Range("I71").Value = 6001
$ws.Range("J71").Value = 8334.333000000001
$ws.Range("K71").Value = 30005
$ws.Range("L71").Value = 41671.665
$ws.Range("M71").Value = -26261
$ws.Range("N71").Value = -49159.665

$ws.Range("H93").Value = 1996
$ws.Range("I93").Value = 1995
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1995
$ws.Range("L93").Value = 2000
$ws.Range("N93").Value = -4496
$ws.Range("M93").Value = -747

$ws.Range("H100").Value = 1349.7778
$ws.Range("I100").Value = 1330
$ws.Range("J100").Value = 1374.5
$ws.Range("K100").Value = 1330
$ws.Range("L100").Value = 1374.5
$ws.Range("M100").Value = -789
$ws.Range("N100").Value = -2456.5

$ws.Range("H101").Value = 22181
$ws.Range("J101").Value = 22181
$ws.Range("L101").Value = 22181
$ws.Range("N101").Value = -28671

$ws.Range("H136").Value = 1465.3334
$ws.Range("I136").Value = 1367.25
$ws.Range("K136").Value = 4101.75
$ws.Range("M136").Value = -1551.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344

$ws.Range("H126").Value = 2413.5
$ws.Range("J126").Value = 1750
$ws.Range("L126").Value = 5250
$ws.Range("N126").Value = -10190
